# Added new references for ATH, maize, rice and tomato on top of dropdown
# (plus a new Ptrichocarpa row appended at the bottom)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows above the current row 2 (pushes existing data rows down to 6..24)
$ws.Rows("2:5").Insert()

# New reference: Arabidopsis / ATH
$ws.Range("A2").Value = "Arabidopsis_1(new).csv"
$ws.Range("B2").Value = "Arabidopsis"
$ws.Range("C2").Value = "?"
$ws.Range("D2").Value = "?"
$ws.Range("E2").Value = "?"
$ws.Range("F2").Value = "?"
$ws.Range("G2").Value = "ATH****"
$ws.Range("H2").Value = "Yes"

# New reference: Zeamays (maize)
$ws.Range("A3").Value = "Zeamays_1(new).csv"
$ws.Range("B3").Value = "Zeamays"
$ws.Range("C3").Value = "?"
$ws.Range("D3").Value = "?"
$ws.Range("E3").Value = "?"
$ws.Range("F3").Value = "?"
$ws.Range("G3").Value = "?"
$ws.Range("H3").Value = "Yes"

# New reference: Oryza (rice)
$ws.Range("A4").Value = "Oryza_1(new).csv"
$ws.Range("B4").Value = "Oryza"
$ws.Range("C4").Value = "?"
$ws.Range("D4").Value = "?"
$ws.Range("E4").Value = "?"
$ws.Range("F4").Value = "?"
$ws.Range("G4").Value = "?"
$ws.Range("H4").Value = "Yes"

# New reference: Solanum_lycopersicum (tomato)
$ws.Range("A5").Value = "Solanum_lycopersicum_1(new).csv"
$ws.Range("B5").Value = "Solanum_lycopersicum"
$ws.Range("C5").Value = "?"
$ws.Range("D5").Value = "?"
$ws.Range("E5").Value = "?"
$ws.Range("F5").Value = "?"
$ws.Range("G5").Value = "?"
$ws.Range("H5").Value = "Yes"

# New Ptrichocarpa row appended at the bottom (row 25);
# B25 keeps the pasted-in Arial / #222222 font styling from the source diff
$ws.Range("A25").Value = "?"
$ws.Range("B25").Value = "Ptrichocarpa"
$ws.Range("C25").Value = "?"
$ws.Range("D25").Value = "?"
$ws.Range("E25").Value = "?"
$ws.Range("F25").Value = "?"
$ws.Range("G25").Value = "Potri.001G134900.1.p"
$ws.Range("H25").Value = "No"
$ws.Range("B25").Font.Name = "Arial"
$ws.Range("B25").Font.Size = 12
$ws.Range("B25").Font.Color = 2236962

# Leave the cursor where the author's last edit landed
$ws.Range("H4").Select()

Write-Output "Applied reference-table updates"
